# burndown chart and classDiagram
#
# 1. Add column headers "ideal" (C1) and "actual" (D1) used as the chart's
#    legend/series labels.
# 2. Reposition the existing burndown chart (move only, size unchanged):
#    from col2/row6 .. col13/row25  ->  col6/row5 .. col17/row24
# 3. Update the active selection to D4 (matches the edited workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- New header labels (creates the shared-strings entries "ideal"/"actual") --
$ws.Range("C1").Value = "ideal"
$ws.Range("D1").Value = "actual"

# -- Move the chart to its new anchor position (same width/height) --
$co = $ws.ChartObjects().Item(1)
$co.Left = 371.025
$co.Top = 75

# -- Match the saved selection/active cell --
$ws.Range("D4").Select()
